$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 11
# from serial 45233 (2023-11-03) to 45243 (2023-11-13).
for ($row = 2; $row -le 11; $row++) {
    $ws.Range("C$row").Value = 45243
}
